$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Hiperlink" cell format applied to B3 (for dorisluvizute@gmail.com) needs to
# move to B2 along with the data, since the rows are being reordered.
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Row 2 <- old row 3 data (Bruno  da Silva)
$ws.Range("A2").Value = "Bruno  da Silva"
$ws.Range("B2").Value = "dorisluvizute@gmail.com"
$ws.Range("C2").Value = "ENVIADO"

# Row 3 <- old row 4 data (Matheus Diniz); column B reverts to the default/Normal style
$ws.Range("A3").Value = "Matheus Diniz"
$ws.Range("B3").Value = "matheusinhodinizinho@gmail.com"
$ws.Range("C3").Value = "ENVIADO"
$ws.Range("B3").Style = "Normal"

# Row 4 <- old row 2 data (Doris Andressa), with the name's accented character fixed
$ws.Range("A4").Value = "Dóris Andressa Moura Luvizute"
$ws.Range("B4").Value = "doriluvizute@gmail.com"
$ws.Range("C4").Value = "Email inválido"

# Update the active selection/cursor position to B4, matching the saved sheet view
$ws.Range("B4").Select()
